# Apply the "Добавить строки 19:15-19:45" update.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column B in this table alternates between two fixed "Посыл / Заповедь" values
# (a plain row and an encrypted/"hash" row). Grab the encrypted value already
# present in B2 so the new rows reuse the exact same shared string.
$hashValue = $ws.Range("B2").Text

# Row 7 switches from the "plain" B-value to the encrypted B-value, and its
# time slot moves from 19:0-19:5 to 19:15-19:20.
$ws.Range("B7").Value = $hashValue
$ws.Range("C7").Value = "19:15-19:20"

# Append five new 5-minute slots (rows 8-12), all using the encrypted B-value.
$ws.Range("B8").Value = $hashValue
$ws.Range("C8").Value = "19:20-19:25"

$ws.Range("B9").Value = $hashValue
$ws.Range("C9").Value = "19:25-19:30"

$ws.Range("B10").Value = $hashValue
$ws.Range("C10").Value = "19:30-19:35"

$ws.Range("B11").Value = $hashValue
$ws.Range("C11").Value = "19:35-19:40"

$ws.Range("B12").Value = $hashValue
$ws.Range("C12").Value = "19:40-19:45"

# Restore the selection left behind in the source file after these edits.
$ws.Range("B21").Select()
